$d = $word.ActiveDocument

$oldSnippet = '<w:pict w14:anchorId="2E114FD8"><v:shapetype id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="_x0000_s2050" type="#_x0000_t202" style="position:absolute;margin-left:249.4pt;margin-top:47pt;width:90.7pt;height:99.2pt;z-index:-1;visibility:visible;mso-wrap-distance-left:9pt;mso-wrap-distance-top:3.6pt;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:3.6pt;mso-position-horizontal-relative:text;mso-position-vertical-relative:text;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:top" stroked="f"><v:textbox><w:txbxContent><w:p w14:paraId="6930CC36" w14:textId="77777777" w:rsidR="008B05DE" w:rsidRDefault="008B05DE" w:rsidP="00EA2A76"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r w:rsidRPr="006F7F3C"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>{%DemisIdQrImage}</w:t></w:r></w:p><w:p w14:paraId="10171D0B" w14:textId="77777777" w:rsidR="008B05DE" w:rsidRPr="000E6D63" w:rsidRDefault="008B05DE" w:rsidP="00EA2A76"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr></w:pPr><w:r w:rsidRPr="000E6D63"><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr><w:t>Meldungs-ID</w:t></w:r></w:p></w:txbxContent></v:textbox></v:shape></w:pict>'
$newSnippet = '<w:pict w14:anchorId="0D08169D"><v:shapetype id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="_x0000_s2052" type="#_x0000_t202" style="position:absolute;margin-left:249.4pt;margin-top:47pt;width:90.7pt;height:99.2pt;z-index:-1;visibility:visible;mso-wrap-distance-left:9pt;mso-wrap-distance-top:3.6pt;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:3.6pt;mso-position-horizontal-relative:text;mso-position-vertical-relative:text;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:top" stroked="f"><v:textbox><w:txbxContent><w:p w14:paraId="6930CC36" w14:textId="77777777" w:rsidR="008B05DE" w:rsidRDefault="008B05DE" w:rsidP="00EA2A76"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>{#</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>DemisIdQrImage</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r><w:r w:rsidRPr="006F7F3C"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>{%DemisIdQrImage}</w:t></w:r></w:p><w:p w14:paraId="10171D0B" w14:textId="77777777" w:rsidR="008B05DE" w:rsidRPr="000E6D63" w:rsidRDefault="008B05DE" w:rsidP="00EA2A76"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr></w:pPr><w:r w:rsidRPr="000E6D63"><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr><w:t>Meldungs-ID</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>DemisIdQrImage</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p></w:txbxContent></v:textbox></v:shape></w:pict>'

$xml = $d.Content.WordOpenXML

if ($xml.IndexOf($oldSnippet) -lt 0) {
    throw "old snippet not found in WordOpenXML"
}

$newXml = $xml.Replace($oldSnippet, $newSnippet)

$d.Content.WordOpenXML = $newXml

Write-Output "Replacement applied."
